$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 263
$ws.Range("I2").Value = 704
$ws.Range("J2").Value = 2975
$ws.Range("K2").Value = 14
$ws.Range("L2").Value = 833
$ws.Range("M2").Value = 50
$ws.Range("N2").Value = 532
$ws.Range("P2").Value = 16
$ws.Range("Q2").Value = 3
$ws.Range("R2").Value = 40
$ws.Range("S2").Value = 288
$ws.Range("T2").Value = 542
$ws.Range("U2").Value = 47
$ws.Range("V2").Value = 4524
$ws.Range("X2").Value = 4548
$ws.Range("Y2").Value = 8
$ws.Range("Z2").Value = 68
$ws.Range("AA2").Value = 33
